# Heropoints.xlsx update — "created personal sites for every profile"
#
# 1) Update hero point totals (column C) for a handful of heroes that
#    previously sat at 0.
# 2) Append two new heroes (Grimstroke, Mars) as new rows at the bottom
#    of the table, each with their own id (column A) and point total
#    (column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated point totals for existing heroes ---------------------------
$ws.Range("C1").Value  = -8     # Anti-Mage
$ws.Range("C2").Value  = -2     # Axe
$ws.Range("C10").Value = -15    # Morphling
$ws.Range("C11").Value = -2     # Shadow Fiend
$ws.Range("C86").Value = 4      # Rubick

# --- New heroes appended to the table ------------------------------------
$ws.Range("A120").Value = 121
$ws.Range("B120").Value = "Grimstroke"
$ws.Range("C120").Value = 5

$ws.Range("A121").Value = 129
$ws.Range("B121").Value = "Mars"
$ws.Range("C121").Value = 4
